$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Tumor" query text (SamplesTab row, column B) so that it references
# samp.sample_tumor_status directly instead of the aliased `tumor` collection value.
$newTumorQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["Childhood Cancer Data Initiative (CCDI): Free the Data: Open Sharing of Comprehensive Genomic Childhood Cancer Datasets (Kansas)"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newTumorQuery

# Move the selection to B4 (FilesTab row), matching the saved selection state.
$ws.Range("B4").Select()
